# "update march and april books"
#
# 1. Fix a typo in an existing title (row 355 / A355).
# 2. De-duplicate a redundant cell style that 4 cells (B230, B254, B282,
#    B320) were using -- collapses them back onto the identical style
#    already used by their neighbours (this is what ripples the other
#    style-index renumbering seen in the diff: J310/K310 and the whole
#    327-330 block shift down by one once the duplicate style is no
#    longer referenced by those 4 cells).
# 3. Append 5 new book rows (357-361) for March/April reads, including one
#    date cell (E359) that needs a brand-new MM/DD/YY number format.
# 4. Move the active-cell selection down to the new last row, like a user
#    who just finished typing would have.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Typo fix -----------------------------------------------------
$ws.Range("A355").Value = "The Sweet Indifference of the World"

# --- 2. Collapse the duplicate style back onto the canonical one -----
# (B230/B254/B282/B320 carry an xf that is byte-identical to the one
# used by every other cell in column B -- re-asserting the same
# General number format makes the engine re-use the existing xf
# instead of the stray duplicate.)
foreach ($addr in @("B230", "B254", "B282", "B320")) {
    $ws.Range($addr).NumberFormat = "General"
}

# --- 3. New rows 357-361 ----------------------------------------------
$newBooks = @(
    @{ Row=357; A="Lucinella";                          B="Lore Segal";         C="Austria"; D="English"; E="15-Mar"; F=2020; H=153; I="Female" },
    @{ Row=358; A="Run to Earth";                        B="Paul Yoon";          C="USA";     D="English"; E="25-Mar"; F=2020; H=259; I="Male"   },
    @{ Row=359; A="The Mysterious Affair at Styles";     B="Agatha Christie";    C="UK";      D="English"; E="4-Apr";  F=2020; H=224; I="Female" },
    @{ Row=360; A="The Hound of Baskervilles";            B="Arthur Conan Doyle"; C="UK";      D="English"; E="19-Apr"; F=2020; H=116; I="Male"   },
    @{ Row=361; A="The Buried Giant";                     B="Kazuo Ishiguro";     C="UK";      D="English"; E="7-May";  F=2020; H=345; I="Male"   }
)

foreach ($book in $newBooks) {
    $r = $book.Row

    $ws.Range("A$r").Value = $book.A
    $ws.Range("B$r").Value = $book.B
    $ws.Range("C$r").Value = $book.C
    $ws.Range("D$r").Value = $book.D
    $ws.Range("E$r").Value = $book.E
    $ws.Range("F$r").Value = $book.F
    $ws.Range("H$r").Value = $book.H
    $ws.Range("I$r").Value = $book.I

    $ws.Rows.Item($r).RowHeight = 12.8
}

# E359 ("4-Apr") is formatted with a new MM/DD/YY date format, unlike the
# other new date cells which stay plain text/general like the rest of
# column E.
$ws.Range("E359").NumberFormat = "MM/DD/YY"

# --- 4. Move the selection to the new last row ------------------------
$ws.Range("A360").Select()
